$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028678020988893
$ws.Range("D2").Value = 1.036796528020931
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.043812063205777
$ws.Range("I2").Value = 1.03514486986706
$ws.Range("J2").Value = 1.033829008540978
$ws.Range("K2").Value = 1.039589244897956
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.046584875411446
$ws.Range("N2").Value = 1.015302514257191

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029397031238671
$ws.Range("D3").Value = 1.037350792930045
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.044542329342464
$ws.Range("I3").Value = 1.03527165712086
$ws.Range("J3").Value = 1.034189910365411
$ws.Range("K3").Value = 1.039953862585193
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.047126463739498
$ws.Range("N3").Value = 1.015422212273532

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029862926031738
$ws.Range("D4").Value = 1.037710030552985
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.045015937048024
$ws.Range("I4").Value = 1.035352819971266
$ws.Range("J4").Value = 1.034423379942194
$ws.Range("K4").Value = 1.040189678436958
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.047477313166697
$ws.Range("N4").Value = 1.0154996327523

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030058941083473
$ws.Range("D5").Value = 1.037861193587417
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.045215297211179
$ws.Range("I5").Value = 1.035386730132299
$ws.Range("J5").Value = 1.034521515340477
$ws.Range("K5").Value = 1.04028878620341
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.047624905516921
$ws.Range("N5").Value = 1.015532172168708

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030091861779084
$ws.Range("D6").Value = 1.037886582650716
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.045248785575238
$ws.Range("I6").Value = 1.03539241142013
$ws.Range("J6").Value = 1.034537991775299
$ws.Range("K6").Value = 1.040305425079259
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.047649692453933
$ws.Range("N6").Value = 1.015537635188483

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029865544596285
$ws.Range("D7").Value = 1.037712049855068
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.045018599907271
$ws.Range("I7").Value = 1.035353273909345
$ws.Range("J7").Value = 1.034424691293616
$ws.Range("K7").Value = 1.040191002836356
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.04747928493095
$ws.Range("N7").Value = 1.015500067578185

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028920878664422
$ws.Range("D8").Value = 1.036983720843203
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.044058635763011
$ws.Range("I8").Value = 1.035187899346225
$ws.Range("J8").Value = 1.033950988226992
$ws.Range("K8").Value = 1.039712492362042
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.046767822640748
$ws.Range("N8").Value = 1.015342973142136

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027261302408365
$ws.Range("D9").Value = 1.035704927516679
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.042375403051982
$ws.Range("I9").Value = 1.034889804844449
$ws.Range("J9").Value = 1.033115876673135
$ws.Range("K9").Value = 1.038868468159297
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.045517322748074
$ws.Range("N9").Value = 1.015065926749121

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026158432421936
$ws.Range("D10").Value = 1.034855620741358
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.041258996110718
$ws.Range("I10").Value = 1.034686626646177
$ws.Range("J10").Value = 1.032558954807629
$ws.Range("K10").Value = 1.038305313036863
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.044685909543512
$ws.Range("N10").Value = 1.014881104643144

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02568173685768
$ws.Range("D11").Value = 1.034488651038548
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.040776971934912
$ws.Range("I11").Value = 1.034597602652468
$ws.Range("J11").Value = 1.032317775325941
$ws.Range("K11").Value = 1.038061366473596
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.04432645633133
$ws.Range("N11").Value = 1.014801050906851

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025504801364837
$ws.Range("D12").Value = 1.034352462264228
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.040598137503626
$ws.Range("I12").Value = 1.034564378716603
$ws.Range("J12").Value = 1.032228187582668
$ws.Range("K12").Value = 1.037970740782662
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.04419302457314
$ws.Range("N12").Value = 1.014771312168619

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025542748691195
$ws.Range("D13").Value = 1.034381669750178
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.040636488518625
$ws.Range("I13").Value = 1.034571512431551
$ws.Range("J13").Value = 1.032247404569216
$ws.Range("K13").Value = 1.037990180869456
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.044221642253509
$ws.Range("N13").Value = 1.0147776913669

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025667108641076
$ws.Range("D14").Value = 1.034477391169797
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.040762185110895
$ws.Range("I14").Value = 1.034594859538903
$ws.Range("J14").Value = 1.032310370027375
$ws.Range("K14").Value = 1.038053875584078
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.044315425072117
$ws.Range("N14").Value = 1.014798592756358

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025743748261168
$ws.Range("D15").Value = 1.034536384273387
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.040839658935031
$ws.Range("I15").Value = 1.034609223752555
$ws.Range("J15").Value = 1.032349164774565
$ws.Range("K15").Value = 1.038093118310616
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.044373219104698
$ws.Range("N15").Value = 1.014811470380146

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026190087317266
$ws.Range("D16").Value = 1.034879992042858
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.041291015900927
$ws.Range("I16").Value = 1.034692512892145
$ws.Range("J16").Value = 1.032574960608043
$ws.Range("K16").Value = 1.038321501053695
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.044709777108527
$ws.Range("N16").Value = 1.014886417069099

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026470294248017
$ws.Range("D17").Value = 1.035095739941536
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.041574513320068
$ws.Range("I17").Value = 1.034744478331331
$ws.Range("J17").Value = 1.032716589660266
$ws.Range("K17").Value = 1.038464734654677
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.044921040773581
$ws.Range("N17").Value = 1.014933422908717

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026633816508238
$ws.Range("D18").Value = 1.035221657714561
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.041740006273276
$ws.Range("I18").Value = 1.034774687870204
$ws.Range("J18").Value = 1.03279919663609
$ws.Range("K18").Value = 1.038548270965607
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.045044320638762
$ws.Range("N18").Value = 1.01496083821991

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026689587268556
$ws.Range("D19").Value = 1.035264605193868
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.04179645770016
$ws.Range("I19").Value = 1.034784971383888
$ws.Range("J19").Value = 1.032827362954652
$ws.Range("K19").Value = 1.038576753046062
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.045086364876859
$ws.Range("N19").Value = 1.014970185711153

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026440222168586
$ws.Range("D20").Value = 1.035072584391634
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.041544082883173
$ws.Range("I20").Value = 1.034738913375251
$ws.Range("J20").Value = 1.032701394488763
$ws.Range("K20").Value = 1.038449368009287
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.044898368650359
$ws.Range("N20").Value = 1.014928379871806

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025630484117118
$ws.Range("D21").Value = 1.03444920026526
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.040725164758991
$ws.Range("I21").Value = 1.034587988710652
$ws.Range("J21").Value = 1.032291828339731
$ws.Range("K21").Value = 1.038035119427493
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.044287806003107
$ws.Range("N21").Value = 1.014792437903117

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025122125531924
$ws.Range("D22").Value = 1.034057949875809
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.040211499281418
$ws.Range("I22").Value = 1.034492191528101
$ws.Range("J22").Value = 1.032034301476486
$ws.Range("K22").Value = 1.037774589967006
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.043904414740392
$ws.Range("N22").Value = 1.014706947416763

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025391543688539
$ws.Range("D23").Value = 1.034265292453596
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.040483686557474
$ws.Range("I23").Value = 1.034543060942002
$ws.Range("J23").Value = 1.032170822438118
$ws.Range("K23").Value = 1.037912708183033
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.044107610268589
$ws.Range("N23").Value = 1.014752269125681

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02645381019284
$ws.Range("D24").Value = 1.035083047154381
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.041557832674944
$ws.Range("I24").Value = 1.034741428251609
$ws.Range("J24").Value = 1.032708260541701
$ws.Range("K24").Value = 1.038456311564209
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.044908613042372
$ws.Range("N24").Value = 1.01493065861059

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027689732159178
$ws.Range("D25").Value = 1.03603496701669
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.042809555862763
$ws.Range("I25").Value = 1.034967656801177
$ws.Range("J25").Value = 1.033331810667395
$ws.Range("K25").Value = 1.039086758048803
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.045840217832307
$ws.Range("N25").Value = 1.01513757383086
